# Update the NATMI LR-pair stats (Grn-Tnfrsf1a) sheet with newly computed
# TPM-based values. Columns G,H,I,J (ligand average/total expression and its
# derived specificities) depend only on the "Sending cluster" (column A).
# Columns M,N,O,P (receptor average/total expression and its derived
# specificities) depend only on the "Target cluster" (column D). Columns
# Q,R,S,T (edge weights/specificities) are simply the products of the
# corresponding ligand/receptor columns: Q=G*M, R=H*N, S=I*O, T=J*P.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ligand-side values (G, H, I, J) keyed by Sending cluster name.
$ligandStats = @{
    "ECs"               = @(32.00264033333334, 96.00792100000001, 0.02419001798940439, 0.02433952891158457)
    "FAPs"              = @(100.282963,         300.848889,        0.07580145430919519, 0.07626995934880827)
    "Inflammatory-Mac"  = @(473.968811,          1421.906433,       0.3582615042098434,  0.360475806319893)
    "MuSCs"             = @(24.3798835,          48.759767,         0.01842816137361988, 0.01236137337687614)
    "Resolving-Mac"     = @(692.3345543333334,   2077.003663,       0.5233188621179371,  0.5265533320428379)
}

# New receptor-side values (M, N, O, P) keyed by Target cluster name.
$receptorStats = @{
    "ECs"               = @(36.95112266666666,  110.853368,  0.1740115908809209,  0.1775751473829744)
    "FAPs"              = @(51.86084433333334,  155.582533,  0.2442250025331967,  0.2492264486514428)
    "Inflammatory-Mac"  = @(64.73785366666665,  194.213561,  0.3048658902295037,  0.3111091917238571)
    "MuSCs"             = @(12.7841595,         25.568319,   0.06020363583370166, 0.04095769119761797)
    "Resolving-Mac"     = @(46.01464833333333,  138.043945,  0.216693880522677,   0.2211315210441077)
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $sendingCluster = $ws.Cells.Item($r, 1).Value2
    $targetCluster  = $ws.Cells.Item($r, 4).Value2

    $ligand = $ligandStats[$sendingCluster]
    $receptor = $receptorStats[$targetCluster]

    $g = $ligand[0]
    $h = $ligand[1]
    $i = $ligand[2]
    $j = $ligand[3]

    $m = $receptor[0]
    $n = $receptor[1]
    $o = $receptor[2]
    $p = $receptor[3]

    $ws.Cells.Item($r, 7).Value  = $g
    $ws.Cells.Item($r, 8).Value  = $h
    $ws.Cells.Item($r, 9).Value  = $i
    $ws.Cells.Item($r, 10).Value = $j

    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p

    $ws.Cells.Item($r, 17).Value = $g * $m
    $ws.Cells.Item($r, 18).Value = $h * $n
    $ws.Cells.Item($r, 19).Value = $i * $o
    $ws.Cells.Item($r, 20).Value = $j * $p
}

Write-Output "Updated $($lastRow - 1) rows with new TPM-derived values"
